$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values, entered in the order that reproduces the shared-string
# --- table layout of the target workbook (column-by-column entry) ---

# A2:A3 -> "Req_ PO3_DGW_CYRS_02_V01"
$ws.Range("A2").Value = "Req_ PO3_DGW_CYRS_02_V01"
$ws.Range("A3").Value = "Req_ PO3_DGW_CYRS_02_V01"

# F2:F5 -> "23/01/2020"
$ws.Range("F2").Value = "23/01/2020"
$ws.Range("F3").Value = "23/01/2020"
$ws.Range("F4").Value = "23/01/2020"
$ws.Range("F5").Value = "23/01/2020"

# B2:B3 -> "Question in the alarm part"
$ws.Range("B2").Value = "Question in the alarm part"
$ws.Range("B3").Value = "Question in the alarm part"

# D2 -> first alarm question
$ws.Range("D2").Value = "what is the time ultil the alarm stops?"

# D3 -> second alarm question
$ws.Range("D3").Value = "what is the tone of the alarm buzzer"

# A4:A5 -> "Req_ PO3_DGW_CYRS_06_V01"
$ws.Range("A4").Value = "Req_ PO3_DGW_CYRS_06_V01"
$ws.Range("A5").Value = "Req_ PO3_DGW_CYRS_06_V01"

# B4:B5 -> "Question regarding the buttons"
$ws.Range("B4").Value = "Question regarding the buttons"
$ws.Range("B5").Value = "Question regarding the buttons"

# D4 -> buttons question 1
$ws.Range("D4").Value = "What is the action that happens when each of the 3 buttons are pressed in each mode?"

# D5 -> buttons question 2
$ws.Range("D5").Value = "Does each click on the buttons makes sound from the buzzer?"

# C2:C5 -> sequence numbers
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 4

# --- Formatting ---------------------------------------------------------
# Build each distinct format on an out-of-the-way scratch cell (column Z,
# far past the used range) and copy/paste the format onto the real cells.
# This keeps every destination range's final look correct while only
# materialising each new style once.

# A2:A5 -> bold, 8pt, Calibri Light, centered horizontally, top-aligned vertically
$ws.Range("Z1").Font.Name = "Calibri Light"
$ws.Range("Z1").Font.Size = 8
$ws.Range("Z1").Font.Bold = $true
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").VerticalAlignment = -4160

# C2:C5 -> centered horizontally (default font)
$ws.Range("Z2").HorizontalAlignment = -4108

# B2:B5, D4 -> 8pt Arial
$ws.Range("Z4").Font.Size = 8

# D2:D3 -> 9pt Arial
$ws.Range("Z5").Font.Size = 9

$ws.Range("Z1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)

$ws.Range("Z2").Copy()
$ws.Range("C2:C5").PasteSpecial(-4122)

$ws.Range("Z4").Copy()
$ws.Range("B2:B5").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("Z5").Copy()
$ws.Range("D2:D3").PasteSpecial(-4122)

# Remove the scratch column entirely so it leaves no trace behind.
$ws.Range("Z1:Z5").EntireColumn.Delete()

# Bookmark-style defined name (sheet-scoped) pointing at the first new row,
# mirroring the Word cross-reference anchor carried over with this content.
$ws.Names.Add("_Toc30617783", "=Sheet1!`$A`$2")

$wb.Save()
